# Add 2022-Q3 sheet + data, matching commit "feat: add 2022-Q3 data"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet right after "总计" and name it "2022-Q3".
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# Fetch the format-template sheet (the existing "2022-Q2" sheet, now shifted
# to position 3) AFTER the insertion so the reference isn't stale.
$template = $wb.Worksheets.Item(3)

# Match page margins used by all the other quarter sheets.
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# Copy header-row styling (bold + border) and the column-A index styling
# from the template sheet so the new sheet matches the others exactly.
$template.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$template.Range("A2").Copy()
$q3.Range("A2:A14").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Fill in the "2022-Q3" fund-holding table.
# ---------------------------------------------------------------------------
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q3.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$rows = @(
    @("009970","财通内需增长12个月定期开放混合","9.38","56.38","2.96","0.2776",4),
    @("501046","财通多策略福鑫定期开放灵活配置混合","2.82","85.55","4.42","0.1246",10),
    @("000017","财通可持续发展主题混合","1.95","90.33","3.56","0.0694",7),
    @("013238","财通均衡一年持有期混合A","1.81","86.47","3.56","0.0644",6),
    @("501026","财通多策略福享混合（LOF）","1.43","88.50","3.56","0.0509",7),
    @("501001","财通多策略精选混合（LOF）","0.78","80.38","3.03","0.0236",10),
    @("006522","财通新兴蓝筹混合A","0.22","94.22","5.47","0.0120",6),
    @("006968","财通行业龙头精选混合C","0.07","90.94","3.63","0.0025",7),
    @("006967","财通行业龙头精选混合A","0.05","90.94","3.63","0.0018",7),
    @("009649","嘉实精选平衡混合A","0.07","58.93","2.29","0.0016",8),
    @("006523","财通新兴蓝筹混合C","0.03","94.22","5.47","0.0016",6),
    @("013239","财通均衡一年持有期混合C","0.04","86.47","3.56","0.0014",6),
    @("009650","嘉实精选平衡混合C","0.05","58.93","2.29","0.0011",8)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $q3.Cells.Item($r, 1).Value = $i

    # B: fund code, C: fund name, D-G: numeric-looking values stored as text
    # (leading "'" forces text storage; Style reset drops the auto-applied
    # "Text" number-format style so the cell matches the plain, style-less
    # text cells used on the other quarter sheets).
    for ($col = 2; $col -le 7; $col++) {
        $cell = $q3.Cells.Item($r, $col)
        $cell.Value = "'" + $row[$col - 2]
        $cell.Style = "Normal"
    }

    # H: real number
    $q3.Cells.Item($r, 8).Value = $row[6]
}

$q3.Range("A1").Select()

# ---------------------------------------------------------------------------
# 3. Update the "总计" summary sheet: insert the 2022-Q3 row at the top of
#    the data (row 2) and shift the rest down.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Extend the column-A index styling down to the new row 5.
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)

$summaryRows = @(
    @("2022-Q3", 13, 0.63),
    @("2022-Q2", 4, 1.16),
    @("2022-Q1", 15, 3.13),
    @("2021-Q4", 13, 4.57)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
}

$total.Range("A1").Select()
